$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (not auto-converted numbers) for the Price column, matching the
# original inlineStr cell type, then restore the default (unstyled) cell format.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.015.08'
$ws.Range("E2").Value = '  -2.60%  '

$ws.Range("D3").Value = '1.860.88'
$ws.Range("E3").Value = '  -2.22%  '

$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").Value = '306.69'
$ws.Range("E5").Value = '  -1.87%  '

$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").Value = '0.5109'
$ws.Range("E7").Value = '  +3.02%  '

$ws.Range("D8").Value = '0.3739'
$ws.Range("E8").Value = '  -0.87%  '

$ws.Range("D9").Value = '0.07115'
$ws.Range("E9").Value = '  -1.86%  '

$ws.Range("D10").Value = '0.8906'
$ws.Range("E10").Value = '  -1.06%  '

$ws.Range("D11").Value = '20.59'
$ws.Range("E11").Value = '  -2.20%  '

$ws.Range("D12").Value = '0.07543'
$ws.Range("E12").Value = '  -1.11%  '

$ws.Range("D13").Value = '1.856.08'
$ws.Range("E13").Value = '  -2.25%  '

$ws.Range("D14").Value = '5.297'
$ws.Range("E14").Value = '  -2.81%  '

$ws.Range("D15").Value = '88.69'
$ws.Range("E15").Value = '  -3.25%  '

$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.32%  '

$ws.Range("D17").Value = '0.000008386'
$ws.Range("E17").Value = '  -3.48%  '

$ws.Range("D18").Value = '14.07'
$ws.Range("E18").Value = '  -3.01%  '

$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.02%  '

$ws.Range("D20").Value = '27.050.04'
$ws.Range("E20").Value = '  -2.65%  '

$ws.Range("D21").Value = '5.062'
$ws.Range("E21").Value = '  -1.85%  '

$ws.Range("D22").Value = '2.095.94'
$ws.Range("E22").Value = '  -3.03%  '

$ws.Range("D23").Value = '10.53'
$ws.Range("E23").Value = '  -2.53%  '

$ws.Range("D24").Value = '6.474'
$ws.Range("E24").Value = '  -1.53%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '149.28'
$ws.Range("E25").Value = '  -2.37%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '1.841'
$ws.Range("E26").Value = '  -0.16%  '

$ws.Range("D27").Value = '17.97'
$ws.Range("E27").Value = '  -1.82%  '

$ws.Range("D28").Value = '2.091'
$ws.Range("E28").Value = '  -5.19%  '

$ws.Range("D29").Value = '112.87'
$ws.Range("E29").Value = '  -1.80%  '

$ws.Range("D30").Value = '4.682'
$ws.Range("E30").Value = '  -3.66%  '

$ws.Range("D31").Value = '4.651'
$ws.Range("E31").Value = '  -2.87%  '

$ws.Range("D32").Value = '0.09043'
$ws.Range("E32").Value = '  +1.23%  '

$ws.Range("D33").Value = '0.05116'
$ws.Range("E33").Value = '  -3.16%  '

$ws.Range("E34").Value = '  -3.79%  '

$ws.Range("D35").Value = '1.156'
$ws.Range("E35").Value = '  -5.92%  '

$ws.Range("D36").Value = '0.7311'
$ws.Range("E36").Value = '  -6.60%  '

$ws.Range("D37").Value = '0.02049'
$ws.Range("E37").Value = '  -1.08%  '

$ws.Range("D38").Value = '2.501'
$ws.Range("E38").Value = '  -5.13%  '

$ws.Range("D39").Value = '3.049'
$ws.Range("E39").Value = '  -0.19%  '

$ws.Range("D40").Value = '1.072'
$ws.Range("E40").Value = '  -1.88%  '

$ws.Range("D41").Value = '0.5326'
$ws.Range("E41").Value = '  -3.07%  '

$ws.Range("D42").Value = '6.592'
$ws.Range("E42").Value = '  -2.04%  '

$ws.Range("D43").Value = '115.03'
$ws.Range("E43").Value = '  +1.12%  '

$ws.Range("D44").Value = '8.315'
$ws.Range("E44").Value = '  -1.43%  '

$ws.Range("D45").Value = '0.1470'
$ws.Range("E45").Value = '  -2.49%  '

$ws.Range("E46").Value = '  +0.02%  '

$ws.Range("D47").Value = '0.4619'
$ws.Range("E47").Value = '  -3.20%  '

$ws.Range("D48").Value = '10.04'
$ws.Range("E48").Value = '  -4.19%  '

$ws.Range("D49").Value = '1.565'
$ws.Range("E49").Value = '  -3.75%  '

$ws.Range("D50").Value = '36.71'
$ws.Range("E50").Value = '  -0.23%  '

$ws.Range("D51").Value = '64.14'
$ws.Range("E51").Value = '  -4.42%  '

$ws.Range("D2:D51").Style = "Normal"
